$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update numeric values (loc x / loc y columns B and C) ---
$ws.Range("B2").Value = 441
$ws.Range("C2").Value = 430

$ws.Range("B3").Value = 477
$ws.Range("C3").Value = 372

$ws.Range("B4").Value = 276
$ws.Range("C4").Value = 266

$ws.Range("B5").Value = 227
$ws.Range("C5").Value = 450

$ws.Range("B6").Value = 430
$ws.Range("C6").Value = 300

$ws.Range("B7").Value = 586
$ws.Range("C7").Value = 260

$ws.Range("B8").Value = 682
$ws.Range("C8").Value = 220

$ws.Range("B9").Value = 443
$ws.Range("C9").Value = 160
$ws.Range("G9").Value = 238

$ws.Range("B10").Value = 511
$ws.Range("C10").Value = 165
$ws.Range("G10").Value = 298

$ws.Range("B11").Value = 405
$ws.Range("C11").Value = 201

$ws.Range("B12").Value = 267
$ws.Range("C12").Value = 156

$ws.Range("B13").Value = 674
$ws.Range("C13").Value = 267

$ws.Range("B14").Value = 257
$ws.Range("C14").Value = 431

$ws.Range("B15").Value = 300
$ws.Range("C15").Value = 270

# --- Row 15's canteen name changed from "North Spine Plaza" to "Mc Donald" ---
$ws.Range("A15").Value = "Mc Donald"

# --- The Address column (H) is no longer populated for rows 6-15 ---
$ws.Range("H6:H15").Clear()

# --- Rows 6-15 no longer carry an explicit (wrapped-text) row height ---
$ws.Rows("6:15").AutoFit()

# --- Rows 2-5 keep an explicit row height, but it shrinks slightly ---
$ws.Rows("2").RowHeight = 62.4
$ws.Rows("3").RowHeight = 62.4
$ws.Rows("4").RowHeight = 124.8
$ws.Rows("5").RowHeight = 78

# --- The saved selection moved to A3 ---
$null = $ws.Range("A3").Select()
